$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($sheet, $addr, $val) {
    $cell = $sheet.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue $ws 'D2' '30.657.35'
Set-TextValue $ws 'E2' '  +0.58%  '
Set-TextValue $ws 'D3' '2.112.47'
Set-TextValue $ws 'E3' '  +0.35%  '
Set-TextValue $ws 'E4' '  +1.13%  '
Set-TextValue $ws 'D5' '351.70'
Set-TextValue $ws 'E5' '  +4.75%  '
Set-TextValue $ws 'D6' '1.012'
Set-TextValue $ws 'E6' '  +0.97%  '
Set-TextValue $ws 'D7' '0.5265'
Set-TextValue $ws 'E7' '  +0.39%  '
Set-TextValue $ws 'D8' '0.4508'
Set-TextValue $ws 'E8' '  -2.09%  '
Set-TextValue $ws 'D9' '53.67'
Set-TextValue $ws 'E9' '  +0.81%  '
Set-TextValue $ws 'D10' '0.09020'
Set-TextValue $ws 'E10' '  +0.61%  '
Set-TextValue $ws 'E11' '  -0.66%  '
Set-TextValue $ws 'D12' '24.52'
Set-TextValue $ws 'E12' '  +0.19%  '
Set-TextValue $ws 'D13' '2.112.50'
Set-TextValue $ws 'E13' '  +0.97%  '
Set-TextValue $ws 'D14' '6.818'
Set-TextValue $ws 'E14' '  +0.07%  '
Set-TextValue $ws 'D15' '8.017'
Set-TextValue $ws 'E15' '  +0.75%  '
Set-TextValue $ws 'D16' '99.75'
Set-TextValue $ws 'E16' '  +3.29%  '
Set-TextValue $ws 'D17' '0.00001174'
Set-TextValue $ws 'E17' '  +3.40%  '
Set-TextValue $ws 'D18' '1.013'
Set-TextValue $ws 'E18' '  +1.04%  '
Set-TextValue $ws 'D19' '0.06708'
Set-TextValue $ws 'E19' '  +0.97%  '
Set-TextValue $ws 'D20' '19.33'
Set-TextValue $ws 'E20' '  +0.07%  '
Set-TextValue $ws 'E21' '  +0.95%  '
Set-TextValue $ws 'D22' '6.284'
Set-TextValue $ws 'E22' '  -0.31%  '
Set-TextValue $ws 'D23' '30.720.01'
Set-TextValue $ws 'E23' '  +0.61%  '
Set-TextValue $ws 'D24' '12.85'
Set-TextValue $ws 'E24' '  +3.94%  '
Set-TextValue $ws 'D25' '2.392'
Set-TextValue $ws 'E25' '  +1.30%  '
Set-TextValue $ws 'D26' '2.358.94'
Set-TextValue $ws 'E26' '  +0.75%  '
Set-TextValue $ws 'E27' '  -0.15%  '
Set-TextValue $ws 'D28' '165.56'
Set-TextValue $ws 'E28' '  +0.96%  '
Set-TextValue $ws 'D29' '2.522'
Set-TextValue $ws 'E29' '  -1.87%  '
Set-TextValue $ws 'D30' '135.77'
Set-TextValue $ws 'E30' '  +2.20%  '
Set-TextValue $ws 'D31' '1.184'
Set-TextValue $ws 'E31' '  -1.95%  '
Set-TextValue $ws 'E32' '  -0.05%  '
Set-TextValue $ws 'D33' '1.634'
Set-TextValue $ws 'E33' '  -3.57%  '
Set-TextValue $ws 'D34' '6.344'
Set-TextValue $ws 'E34' '  +2.72%  '
Set-TextValue $ws 'D35' '4.015'
Set-TextValue $ws 'E35' '  +2.25%  '
Set-TextValue $ws 'D36' '5.886'
Set-TextValue $ws 'E36' '  +5.90%  '
Set-TextValue $ws 'D37' '10.19'
Set-TextValue $ws 'E37' '  -2.53%  '
Set-TextValue $ws 'D38' '0.02647'
Set-TextValue $ws 'E38' '  +2.87%  '
Set-TextValue $ws 'D39' '0.06823'
Set-TextValue $ws 'E39' '  -0.15%  '
Set-TextValue $ws 'D40' '0.2309'
Set-TextValue $ws 'E40' '  +0.74%  '
Set-TextValue $ws 'D41' '12.54'
Set-TextValue $ws 'E41' '  -2.38%  '
Set-TextValue $ws 'D42' '0.6886'
Set-TextValue $ws 'E42' '  +0.02%  '
Set-TextValue $ws 'D43' '1.281'
Set-TextValue $ws 'E43' '  +2.73%  '
Set-TextValue $ws 'D44' '14.77'
Set-TextValue $ws 'E44' '  +5.04%  '
Set-TextValue $ws 'D45' '2.321'
Set-TextValue $ws 'E45' '  -1.38%  '
Set-TextValue $ws 'D46' '0.6421'
Set-TextValue $ws 'E46' '  +0.52%  '
Set-TextValue $ws 'D47' '3.773'
Set-TextValue $ws 'E47' '  +2.82%  '
Set-TextValue $ws 'D48' '0.00000000354'
Set-TextValue $ws 'E48' '  -0.02%  '
Set-TextValue $ws 'D49' '1.247'
Set-TextValue $ws 'E49' '  -0.11%  '
Set-TextValue $ws 'D50' '0.07277'
Set-TextValue $ws 'E50' '  +2.07%  '
Set-TextValue $ws 'D51' '82.42'
Set-TextValue $ws 'E51' '  -1.47%  '
